# TutorialStartX, TutorialStartZ 글로벌컨스턴트플롯테이블에 추가 (임시값 -1.5, -3)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalConstantFloatTable")

# Make this the active sheet/tab (matches the tabSelected move from sheet1 to sheet2)
$ws.Activate()

$ws.Cells.Item(15, 1).Value = "TutorialStartX"
$ws.Cells.Item(15, 2).Value = -1.5

$ws.Cells.Item(16, 1).Value = "TutorialStartZ"
$ws.Cells.Item(16, 2).Value = -3

# Leave selection on the row below the newly added data
$ws.Range("A17").Select()
